$d = $word.ActiveDocument

# 1) Version history table: date cell '09/09/15' -> 'dd/mm/aa'
$d.Content.Find.Execute('09/09/15', $true, $false, $false, $false, $false, $true, 1, $false, 'dd/mm/aa', 2) | Out-Null

# 2) Objectives paragraph 1 (rewritten)
$d.Content.Find.Execute('Temos como principal objetivo do projeto o ensino da lógica de programação tendo como propósito que 50% do público infantil tenha acesso ao jogo e que o jogo proporcione uma melhora de aprendizagem em 80% deste público.', $true, $false, $false, $false, $false, $true, 1, $false, 'Temos como principal objetivo do projeto o desenvolvimento do raciocínio lógico como auxilio para facilidade de aprendizagem em  lógica de programação tendo como propósito que o jogo proporcione uma melhora na aprendizagem de 70% dos usuarios.', 2) | Out-Null

# 3) Objectives paragraph 2 (rewritten)
$d.Content.Find.Execute('É propósito que no mínimo 70% do público alcançado nos de um retorno positivo de forma que  permaneça jogando por um prazo de 1 ano. Ao final de 1 ano do jogo no mercado temos como objetivo que 80% dos usuários cheguem a fase final. ', $true, $false, $false, $false, $false, $true, 1, $false, 'É propósito que no mínimo 70% do público alcançado nos de um retorno positivo de forma que  permaneça jogando por um prazo de 1 mes. Ao final de 1 mes de jogo no mercado temos como objetivo que 80% dos usuários cheguem a fase final. ', 2) | Out-Null

# 4) Objectives paragraph 3 removed entirely
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.TrimEnd() -eq 'É também tido como objetivo que a procura por cursos de ensino de lógica de programação tenha um aumento de 20% nos próximos 5 anos  e que o índice de reprovação em matérias de calculo seja reduzido em cerca de 50%'.TrimEnd()) {
        $p.Range.Delete()
        break
    }
}

# 5) Scope paragraph (rewritten)
$d.Content.Find.Execute('O projeto trata se da criação de um jogo infantil de caça ao tesouro com conceitos de lógica de programação que será disponibilizado para o público de celulares com a plataforma android.  ', $true, $false, $false, $false, $false, $true, 1, $false, 'O projeto trata se da criação de um jogo infantil de caça ao tesouro com conceitos de lógica de programação, que será disponibilizado para o público infantil, que tenham celulares com a plataforma android.  ', 2) | Out-Null

# 6) Heading 'Visão da Situação Proposta' - strip stale lastRenderedPageBreak by re-writing its own text
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Style.NameLocal -eq 'Heading 2' -and $p.Range.Text.TrimEnd() -eq 'Visão da Situação Proposta') {
        $p.Range.Find.Execute('Visão da Situação Proposta', $true, $false, $false, $false, $false, $true, 1, $false, 'Visão da Situação Proposta', 2) | Out-Null
        break
    }
}

# 7) 'A solução proposta...' paragraph rewritten
$d.Content.Find.Execute('A solução proposta para o crescimento de enisino da lógica de programação para o público infantil é dada através de um jogo de caça ao tesouro, o jogo que será desenvolvido para plataforma android contará com alcance de ... crianças que possuem a plataforma atualmente. O usuário poderá obter conhecimento de lógica de programação de uma forma divertida e sem cobranças a forma de pontuação é um incentivo para o interesse pelo jogo seja cada vez maior. O jogo contará com níveis de dificuldade para que os usuários possam se desenvolver gradativamente e que sem que percebam estejam inseridos nesta forma de linguagem. Atualmente existem diversas tentativas e metodologias para alcançar o público proposto mas não encontramos nenhum que tenha sido desenvolvido nesta plataforma e da forma como proposto acima.', $true, $false, $false, $false, $false, $true, 1, $false, 'A solução proposta para o crescimento de incentivo ao enisino da lógica de programação para o público infantil é dada através de um jogo de caça ao tesouro, o jogo  será desenvolvido para plataforma android contará com alcance de crianças que possuem a plataforma atualmente. O usuário poderá desenvolver o raciocínio lógico de uma forma divertida e sem cobranças, a forma de pontuação é um incentivo para que o interesse pelo jogo seja cada vez maior. O jogo contará com níveis de dificuldade para que os usuários possam se desenvolver gradativamente e que sem que percebam estejam inseridos nesta forma de linguagem. Atualmente existem diversas tentativas e metodologias para alcançar o público proposto, mas não encontramos nenhuma que tenha sido desenvolvido nesta plataforma e da forma como proposto acima.', 2) | Out-Null

# 8) Heading 'Recursos e Prazos' - drop stale lastRenderedPageBreak (moved to 'Usabilidade' in real layout)
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Style.NameLocal -eq 'Heading 3' -and $p.Range.Text.TrimEnd() -eq 'Recursos e Prazos') {
        $p.Range.Find.Execute('Recursos e Prazos', $true, $false, $false, $false, $false, $true, 1, $false, 'Recursos e Prazos', 2) | Out-Null
        break
    }
}

# 9) Usability paragraph - '10 segundos' -> '3 segundos'
$d.Content.Find.Execute('O sistema terá um tempo online de resposta para cada comando dado pelo usuário, o placar também deve ser atualizado simultaneamente a pontuação adquirida. O tempo de transição de uma fase para outra deverá ser de no máximo 10 segundos. ', $true, $false, $false, $false, $false, $true, 1, $false, 'O sistema terá um tempo online de resposta para cada comando dado pelo usuário, o placar também deve ser atualizado simultaneamente a pontuação adquirida. O tempo de transição de uma fase para outra deverá ser de no máximo 3 segundos. ', 2) | Out-Null

# 10) Table cell fix: 'Os personagem' -> 'O personagem'
$d.Content.Find.Execute('Os personagem movimenta-se somente na vertical e horizontal.', $true, $false, $false, $false, $false, $true, 1, $false, 'O personagem movimenta-se somente na vertical e horizontal.', 2) | Out-Null

# 11) Table Grid style: add explicit 10pt (sz/szCs 20) run size
$tg = $d.Styles('TableGrid')
$tg.Font.Size = 10
$tg.Font.SizeBi = 10

Write-Output "edit complete"
